# Auto-generated edit script: updates crypto price/volume table
# to match the commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.140.40"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.805.16"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'311.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5116"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.18%  "
$ws.Range("D8").Value = "'0.3916"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "'0.07791"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.28%  "
$ws.Range("D10").Value = "'1.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'41.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").Value = "'6.361"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "'1.003"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "'20.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "'7.324"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "1.808.60"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "'92.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "'0.00001074"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("D19").Value = "'0.06583"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'17.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "'6.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "28.239.05"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "'11.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'160.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'2.455"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.016.75"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").Value = "'127.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("D31").Value = "'0.1094"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").Value = "'1.054"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.15%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.651"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.553"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.17%  "
$ws.Range("D35").Value = "'0.07053"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("D36").Value = "'9.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "'0.02346"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'0.2169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").Value = "'11.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.52%  "
$ws.Range("D40").Value = "'5.011"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").Value = "'0.6171"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'1.156"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").Value = "'13.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").Value = "'0.5925"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("E46").Value = "  -5.58%  "
$ws.Range("D47").Value = "'3.726"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").Value = "'124.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").Value = "'1.207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").Value = "'1.918"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("D51").Value = "'0.06779"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.79%  "
